# "added support for col/row"
#
# The sheet used to store the source well as a bare numeric "Well Id" in
# column B. This edit switches row 2 over to the new "Source Col/Row"
# column (C) instead, and rolls the "Destination Plate Barcode" column
# (E2:E9) forward from the "jul15" batch to the new "jul16" batch.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: clear the old numeric "Source Well Id" cell and populate the
# new "Source Col/Row" cell instead.
$ws.Range("B2").ClearContents() | Out-Null
$ws.Range("C2").Value = "A2"

# All data rows: Destination Plate Barcode moves to the new batch.
$ws.Range("E2:E9").Value = "ssdest000000141jul16"

# Leave the selection on the cell that was actually edited.
$ws.Range("B2").Select() | Out-Null
